# Add a new "2023" column (T) to the malaria-incidence table, mirroring the
# formatting of the existing "2022" column (S), and tighten columns A:C.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- New column T: copy formatting from column S row by row, then set values ---

$ws.Range("S3").Copy()
$ws.Range("T3").PasteSpecial(-4122)
$ws.Range("T3").Value = 2023

$ws.Range("S4").Copy()
$ws.Range("T4").PasteSpecial(-4122)
$ws.Range("T4").Value = 0

$ws.Range("S5").Copy()
$ws.Range("T5").PasteSpecial(-4122)
$ws.Range("T5").Value = "-"

$ws.Range("S6").Copy()
$ws.Range("T6").PasteSpecial(-4122)
$ws.Range("T6").Value = "-"

$ws.Range("S7").Copy()
$ws.Range("T7").PasteSpecial(-4122)
$ws.Range("T7").Value = "-"

$ws.Range("S8").Copy()
$ws.Range("T8").PasteSpecial(-4122)
$ws.Range("T8").Value = "-"

$ws.Range("S9").Copy()
$ws.Range("T9").PasteSpecial(-4122)
$ws.Range("T9").Value = "-"

$ws.Range("S10").Copy()
$ws.Range("T10").PasteSpecial(-4122)
$ws.Range("T10").Value = "-"

$ws.Range("S11").Copy()
$ws.Range("T11").PasteSpecial(-4122)
$ws.Range("T11").Value = "-"

$ws.Range("S12").Copy()
$ws.Range("T12").PasteSpecial(-4122)
$ws.Range("T12").Value = 0.001731197036190674

$ws.Range("S13").Copy()
$ws.Range("T13").PasteSpecial(-4122)
$ws.Range("T13").Value = "-"

# --- Narrow columns A:C slightly ---
$ws.Range("A1:C1").ColumnWidth = 33.3
